$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 0.003873666666666666
$ws.Cells.Item(2, 8).Value = 0.011621
$ws.Cells.Item(2, 9).Value = 0.000207097047037761
$ws.Cells.Item(2, 10).Value = 0.0002070970470377611
$ws.Cells.Item(2, 13).Value = 2.133443333333334
$ws.Cells.Item(2, 14).Value = 6.40033
$ws.Cells.Item(2, 15).Value = 0.2605947899689859
$ws.Cells.Item(2, 16).Value = 0.2605947899689859
$ws.Cells.Item(2, 17).Value = 0.008264248325555557
$ws.Cells.Item(2, 18).Value = 0.07437823493
$ws.Cells.Item(2, 19).Value = 0.00005396841147600254
$ws.Cells.Item(2, 20).Value = 0.00005396841147600255

# Row 3
$ws.Cells.Item(3, 7).Value = 0.003873666666666666
$ws.Cells.Item(3, 8).Value = 0.011621
$ws.Cells.Item(3, 9).Value = 0.000207097047037761
$ws.Cells.Item(3, 10).Value = 0.0002070970470377611
$ws.Cells.Item(3, 15).Value = 0.5209338844846115
$ws.Cells.Item(3, 16).Value = 0.5209338844846116
$ws.Cells.Item(3, 17).Value = 0.01652038777555555
$ws.Cells.Item(3, 18).Value = 0.14868348998
$ws.Cells.Item(3, 19).Value = 0.0001078838691786732
$ws.Cells.Item(3, 20).Value = 0.0001078838691786732

# Row 4
$ws.Cells.Item(4, 7).Value = 0.003873666666666666
$ws.Cells.Item(4, 8).Value = 0.011621
$ws.Cells.Item(4, 9).Value = 0.000207097047037761
$ws.Cells.Item(4, 10).Value = 0.0002070970470377611
$ws.Cells.Item(4, 13).Value = 1.788586
$ws.Cells.Item(4, 14).Value = 5.365758
$ws.Cells.Item(4, 15).Value = 0.2184713255464024
$ws.Cells.Item(4, 16).Value = 0.2184713255464024
$ws.Cells.Item(4, 17).Value = 0.006928385968666666
$ws.Cells.Item(4, 18).Value = 0.06235547371799999
$ws.Cells.Item(4, 19).Value = 0.0000452447663830853
$ws.Cells.Item(4, 20).Value = 0.00004524476638308532

# Row 5
$ws.Cells.Item(5, 9).Value = 0.9943075488985426
$ws.Cells.Item(5, 10).Value = 0.9943075488985427
$ws.Cells.Item(5, 13).Value = 2.133443333333334
$ws.Cells.Item(5, 14).Value = 6.40033
$ws.Cells.Item(5, 15).Value = 0.2605947899689859
$ws.Cells.Item(5, 16).Value = 0.2605947899689859
$ws.Cells.Item(5, 17).Value = 39.67803797112445
$ws.Cells.Item(5, 18).Value = 357.10234174012
$ws.Cells.Item(5, 19).Value = 0.2591113668697929
$ws.Cells.Item(5, 20).Value = 0.259111366869793

# Row 6
$ws.Cells.Item(6, 9).Value = 0.9943075488985426
$ws.Cells.Item(6, 10).Value = 0.9943075488985427
$ws.Cells.Item(6, 15).Value = 0.5209338844846115
$ws.Cells.Item(6, 16).Value = 0.5209338844846116
$ws.Cells.Item(6, 19).Value = 0.5179684938200906
$ws.Cells.Item(6, 20).Value = 0.5179684938200908

# Row 7
$ws.Cells.Item(7, 9).Value = 0.9943075488985426
$ws.Cells.Item(7, 10).Value = 0.9943075488985427
$ws.Cells.Item(7, 13).Value = 1.788586
$ws.Cells.Item(7, 14).Value = 5.365758
$ws.Cells.Item(7, 15).Value = 0.2184713255464024
$ws.Cells.Item(7, 16).Value = 0.2184713255464024
$ws.Cells.Item(7, 17).Value = 33.26433944310133
$ws.Cells.Item(7, 18).Value = 299.379054987912
$ws.Cells.Item(7, 19).Value = 0.2172276882086589
$ws.Cells.Item(7, 20).Value = 0.217227688208659

# Row 8
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.1026013333333333
$ws.Cells.Item(8, 8).Value = 0.307804
$ws.Cells.Item(8, 9).Value = 0.005485354054419671
$ws.Cells.Item(8, 10).Value = 0.005485354054419673
$ws.Cells.Item(8, 13).Value = 2.133443333333334
$ws.Cells.Item(8, 14).Value = 6.40033
$ws.Cells.Item(8, 15).Value = 0.2605947899689859
$ws.Cells.Item(8, 16).Value = 0.2605947899689859
$ws.Cells.Item(8, 17).Value = 0.2188941305911112
$ws.Cells.Item(8, 18).Value = 1.97004717532
$ws.Cells.Item(8, 19).Value = 0.00142945468771702
$ws.Cells.Item(8, 20).Value = 0.00142945468771702

# Row 9
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.1026013333333333
$ws.Cells.Item(9, 8).Value = 0.307804
$ws.Cells.Item(9, 9).Value = 0.005485354054419671
$ws.Cells.Item(9, 10).Value = 0.005485354054419673
$ws.Cells.Item(9, 15).Value = 0.5209338844846115
$ws.Cells.Item(9, 16).Value = 0.5209338844846116
$ws.Cells.Item(9, 17).Value = 0.4375734823911112
$ws.Cells.Item(9, 18).Value = 3.93816134152
$ws.Cells.Item(9, 19).Value = 0.002857506795342253
$ws.Cells.Item(9, 20).Value = 0.002857506795342254

# Row 10
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.1026013333333333
$ws.Cells.Item(10, 8).Value = 0.307804
$ws.Cells.Item(10, 9).Value = 0.005485354054419671
$ws.Cells.Item(10, 10).Value = 0.005485354054419673
$ws.Cells.Item(10, 13).Value = 1.788586
$ws.Cells.Item(10, 14).Value = 5.365758
$ws.Cells.Item(10, 15).Value = 0.2184713255464024
$ws.Cells.Item(10, 16).Value = 0.2184713255464024
$ws.Cells.Item(10, 17).Value = 0.1835113083813333
$ws.Cells.Item(10, 18).Value = 1.651601775432
$ws.Cells.Item(10, 19).Value = 0.001198392571360398
$ws.Cells.Item(10, 20).Value = 0.001198392571360399
